$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45192
# (2023-09-23) to 45202 (2023-10-03) for every data row (rows 2 through 103).
for ($row = 2; $row -le 103; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45192) {
        $cell.Value = 45202
    }
}
